$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the header style (bold, centered, bordered)
# from the neighboring G1 header cell so the new column matches formatting.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column (H2:H20) with the per-row flag values.
$saveValues = @(0,0,0,0,0,1,1,0,1,1,0,0,0,1,0,0,0,0,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
